# Append new bank activation codes to the bottom of column A.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @("88fc95", "fb671a", "0239fd", "0a9813", "131a0b")

# Find the first empty row after the existing data in column A.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($ws.Cells.Item(1, 1).Value -eq $null) {
    $lastRow = 0
}

foreach ($value in $values) {
    $lastRow = $lastRow + 1
    $ws.Cells.Item($lastRow, 1).Value = $value
}
